# Spark-Wiring "Core SW" pin mapping — adds a new column M ("Core SW") next
# to the existing column L (renamed from "Core" to "Core HW"), mirroring the
# commit "Added Core SW pin mapping (from Spark-Wiring)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Core SW" column header first (so it lands earlier in the shared
# string table), then rename the existing "Core" header to "Core HW".
$ws.Range("M1").Value = "Core SW"
$ws.Range("L1").Value = "Core HW"

# New "Core SW" pin numbers for column M.
$ws.Range("M2").Value = 10
$ws.Range("M3").Value = 11
$ws.Range("M4").Value = 19
$ws.Range("M5").Value = 18
$ws.Range("M6").Value = 12
$ws.Range("M7").Value = 13
$ws.Range("M8").Value = 14
$ws.Range("M9").Value = 15
$ws.Range("M10").Value = 8
$ws.Range("M11").Value = 9
$ws.Range("M12").Value = 20
$ws.Range("M15").Value = 7
$ws.Range("M16").Value = 6
$ws.Range("M17").Value = 5
$ws.Range("M18").Value = 16
$ws.Range("M19").Value = 17
$ws.Range("M21").Value = 4
$ws.Range("M22").Value = 3
$ws.Range("M23").Value = 2
$ws.Range("M24").Value = 1
$ws.Range("M25").Value = 0

# Column L no longer needs to be auto-fit-wide now that it only holds short
# "Core HW" labels; narrow it back down to (near) the sheet's default width.
$ws.Columns.Item(12).ColumnWidth = 8.5

# Re-apply the AutoFilter over the now-wider A1:M49 range, and update the
# hidden _FilterDatabase defined name to match.
$ws.AutoFilterMode = $false
$ws.Range("A1:M49").AutoFilter() | Out-Null
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$M`$49"

# Reset the window scroll position (drop the stale topLeftCell="A2") and
# move the active selection to O31, matching the saved view state.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("O31").Select() | Out-Null
